# Update "grille_evaluation" sheet:
#  - Fill in the evaluation grade for the two new criteria rows (13 and 14)
#    with "Tout est réussi", which drives the VLOOKUP results in D13/D14.
#  - Set the student name (B15) to "Tremblay, Charles-Étienne".
#  - Move the active selection to C4 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grille_evaluation")

$ws.Range("C13").Value = "Tout est réussi"
$ws.Range("C14").Value = "Tout est réussi"

$ws.Range("B15").Value = "Tremblay, Charles-Étienne"

$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
